# Applies the commit "Added RepositoryAbstractPage class, test for checking issues creating"
# to the TestCases.xlsx workbook.

$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("TestCases")
$wsData  = $wb.Worksheets.Item("TestData")

# --- TestCases sheet: rename existing test cases ---
$wsCases.Range("A11").Value = "4. testAddRepository"
$wsCases.Range("A16").Value = "5. testDeleteRepository"

# --- TestCases sheet: new test case "6. testAddIssue" fills previously empty rows 19-24 ---
# (cell write order matches the original author's edit order so new shared-string
#  entries land in the same sequence as the recorded workbook)
$wsCases.Range("A19").Value = "6. testAddIssue"
$wsCases.Range("B20").Value = "2. Click on Issues link, check that all sections and welcome message present"
$wsCases.Range("B21").Value = "3. Click the link to create issue, check that Title, Comments fields and Labels, Milestone, Assignee links present"
$wsCases.Range("B22").Value = "4. Fill all fields and confirm creating, check that issue submitted"
$wsCases.Range("B23").Value = "5. Navigate to Issues Section page and check that new issue appeared in the list of issues"
$wsCases.Range("B24").Value = "6. Delete repository"
$wsCases.Range("B19").Value = "1. Log in, add new repository"

# --- selections / active sheet swap: TestCases becomes the active tab ---
$wsData.Range("D5").Select()

$wsCases.Activate()
$excel.ActiveWindow.ScrollRow = 3
$wsCases.Range("E27").Select()
